# Q2E_Peru.xlsx: drop the blank spacer column (old "_1" header in column C)
# and tidy the two stats headers that picked up a duplicated word:
#   "Peru price price"  -> "Peru_priceprice"
#   "Peru points points" -> "Peru_pointspoints"
# Deleting column C shifts the old column D ("points" stats) left into C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty spacer column C; column D (points stats) shifts to C.
$ws.Columns("C").Delete()

# Fix up the header text for the two remaining stat columns.
$ws.Range("B1").Value = "Peru_priceprice"
$ws.Range("C1").Value = "Peru_pointspoints"
